$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Update "driver_Post_ValidationData" test data first so the new
#    shared strings ("400706", "7755668283", "driver.test_83@gmail.com")
#    get registered in that exact order.
# ------------------------------------------------------------------
$validation = $wb.Worksheets.Item("driver_Post_ValidationData")
$validation.Range("B9").Value = "400706"
$validation.Range("B6").Value = "7755668283"
$validation.Range("B5").Value = "driver.test_83@gmail.com"

# ------------------------------------------------------------------
# 2. Update "driver_Post" sheet (email + mobile number test data).
# ------------------------------------------------------------------
$driverPost = $wb.Worksheets.Item("driver_Post")
$driverPost.Range("B5").Value = "7755668283"
$driverPost.Range("B4").Value = "driver.test_83@gmail.com"

# ------------------------------------------------------------------
# 3. Update "driver_verify_Post" sheet (username/email test data).
# ------------------------------------------------------------------
$driverVerify = $wb.Worksheets.Item("driver_verify_Post")
$driverVerify.Range("B2").Value = "driver.test_83@gmail.com"

# ------------------------------------------------------------------
# 4. Insert a brand-new worksheet "driver_LocPinCode" right after
#    "driver_verify_Post" (i.e. before "driver_Post_ValidationData"),
#    for the PATCH request driver location pincode test data.
# ------------------------------------------------------------------
$after = $wb.Worksheets.Item("driver_verify_Post")
$newSheet = $wb.Worksheets.Add($null, $after)
$newSheet.Name = "driver_LocPinCode"

# Copy the cell formatting of a similarly-shaped 2-row Key/Value sheet
# so the new sheet's styles match (header style, text-format cells, etc).
$driverVerify.Range("A1:B3").Copy()
$newSheet.Range("A1:B3").PasteSpecial(-4122)  # xlPasteFormats

$newSheet.Range("A1").Value = "Key"
$newSheet.Range("B1").Value = "Values"
$newSheet.Range("A2").Value = "locationPinCode"
$newSheet.Range("B2").Value = "400706"
$newSheet.Range("A3").Value = "vehicleAvailablity"
$newSheet.Range("B3").Value = "null"

# Hyperlinks.Add with a TextToDisplay argument overwrites the cell's
# displayed text, so re-apply the cell's real value afterwards.
$newSheet.Hyperlinks.Add($newSheet.Range("B3"), "mailto:Admin@demo", [System.Type]::Missing, [System.Type]::Missing, "Admin@demo")
$newSheet.Range("B3").Value = "null"

$newSheet.Range("B2").Select()

# ------------------------------------------------------------------
# 5. Make "driver_verify_Post" the active tab (matches the workbook's
#    new activeTab setting).
# ------------------------------------------------------------------
$driverVerify.Activate()
$driverVerify.Range("B2").Select()

Write-Output "done"
